# payments.xlsx — add payment 71717173 (Cash) 2025-08-20T08:22:30
# and fix the mis-typed phone number that was left as text on row 68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 68: "71717170" was stored as text, should be numeric ---
$ws.Cells.Item(68, 1).Value = 71717170

# --- Append new row 69 for payment 71717173 ---
$row = 69

# Phone number keeps its text representation (matches source data for this row)
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "71717173"

# "amount" column (B) is left blank, like the other recent rows
$ws.Cells.Item($row, 2).Value = ""

$ws.Cells.Item($row, 3).Value = "Cash"
$ws.Cells.Item($row, 4).Value = "2025-08-20T08:22:30"
$ws.Cells.Item($row, 5).Value = 150

# "discount_applied" column (F) is left blank too
$ws.Cells.Item($row, 6).Value = ""

$ws.Cells.Item($row, 7).Value = 127.5
$ws.Cells.Item($row, 8).Value = 22.5
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
